$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns are plain text in the source
# sheet - values such as "29.835.78", "0.570" or "  +2.02%  " are not
# valid Excel numbers, but some of them (e.g. "215.36", "0.570") look
# numeric enough that a plain .Value assignment would get silently
# coerced to a Double - which also eats significant trailing zeros
# (e.g. "0.570" -> 0.57). Force the whole data range to a Text number
# format first so every write below lands as text, then restore the
# Normal style afterwards so no stray formatting is left on cells that
# originally had none (data rows 2-51 have no explicit style).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.835.78"
$ws.Range("E2").Value = "  -0.19%  "

$ws.Range("D3").Value = "1.639.37"
$ws.Range("E3").Value = "  +0.85%  "

$ws.Range("E4").Value = "  +0.26%  "

$ws.Range("D5").Value = "215.36"
$ws.Range("E5").Value = "  +0.36%  "

$ws.Range("E6").Value = "  -0.46%  "

$ws.Range("E7").Value = "  +0.29%  "

$ws.Range("D8").Value = "28.92"
$ws.Range("E8").Value = "  -2.96%  "

$ws.Range("E9").Value = "  +0.70%  "

$ws.Range("D10").Value = "0.0609"
$ws.Range("E10").Value = "  -0.36%  "

$ws.Range("E11").Value = "  -1.59%  "

$ws.Range("D12").Value = "1.874.44"
$ws.Range("E12").Value = "  +0.84%  "

$ws.Range("D13").Value = "1.637.46"
$ws.Range("E13").Value = "  +0.76%  "

$ws.Range("D14").Value = "0.592"
$ws.Range("E14").Value = "  +3.48%  "

$ws.Range("E15").Value = "  +8.01%  "

$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("D17").Value = "29.836.63"
$ws.Range("E17").Value = "  -0.37%  "

$ws.Range("D18").Value = "64.33"
$ws.Range("E18").Value = "  -0.51%  "

$ws.Range("D19").Value = "237.67"
$ws.Range("E19").Value = "  -2.64%  "

$ws.Range("E20").Value = "  -0.49%  "

$ws.Range("E21").Value = "  +0.22%  "

$ws.Range("D22").Value = "9.91"
$ws.Range("E22").Value = "  +2.98%  "

$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("E24").Value = "  +2.21%  "

$ws.Range("D25").Value = "157.34"
$ws.Range("E25").Value = "  -0.08%  "

$ws.Range("D26").Value = "15.59"
$ws.Range("E26").Value = "  -0.50%  "

$ws.Range("E27").Value = "  -1.29%  "

$ws.Range("D28").Value = "6.65"
$ws.Range("E28").Value = "  +0.68%  "

$ws.Range("E29").Value = "  +0.18%  "

$ws.Range("D30").Value = "0.0496"
$ws.Range("E30").Value = "  +1.22%  "

$ws.Range("E31").Value = "  -1.13%  "

$ws.Range("E32").Value = "  +0.88%  "

$ws.Range("E33").Value = "  -1.17%  "

$ws.Range("D34").Value = "1.420.73"
$ws.Range("E34").Value = "  -0.36%  "

$ws.Range("E35").Value = "  +2.52%  "

$ws.Range("D36").Value = "1.02"
$ws.Range("E36").Value = "  -1.36%  "

$ws.Range("E37").Value = "  +1.93%  "

$ws.Range("E38").Value = "  -6.90%  "

$ws.Range("D39").Value = "2.29"
$ws.Range("E39").Value = "  +0.04%  "

# Row 40/41 swap places: Aave and ImmutableX trade ranks (plus updated
# price/volume figures).
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "0.570"
$ws.Range("E40").Value = "  +2.02%  "

$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "76.52"
$ws.Range("E41").Value = "  +10.43%  "

$ws.Range("E42").Value = "  +0.16%  "

$ws.Range("D43").Value = "0.834"
$ws.Range("E43").Value = "  +0.07%  "

$ws.Range("E44").Value = "  -2.46%  "

$ws.Range("E45").Value = "  +0.27%  "

$ws.Range("E46").Value = "  -2.18%  "

$ws.Range("D47").Value = "50.32"
$ws.Range("E47").Value = "  -8.26%  "

$ws.Range("D48").Value = "1.782.48"
$ws.Range("E48").Value = "  +0.88%  "

$ws.Range("E49").Value = "  -1.47%  "

$ws.Range("D50").Value = "93.84"
$ws.Range("E50").Value = "  +5.63%  "

$ws.Range("E51").Value = "  +1.00%  "

# Clear the temporary Text format back to Normal so the saved file has
# no extra styling versus the original (unstyled) data cells.
$ws.Range("D2:E51").Style = "Normal"
